$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 6534
$ws.Range("L3").Value = 7049
$ws.Range("L4").Value = 1757
$ws.Range("L5").Value = 418
$ws.Range("L6").Value = 5779
$ws.Range("L7").Value = 21537
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 435
$ws.Range("L3").Value = 500
$ws.Range("L7").Value = 1420
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L6").Value = 272
$ws.Range("L7").Value = 966
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 107
$ws.Range("L3").Value = 101
$ws.Range("L7").Value = 300
$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 165
$ws.Range("L7").Value = 431
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 114
$ws.Range("L3").Value = 149
$ws.Range("L7").Value = 373
$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L6").Value = 39
$ws.Range("L7").Value = 95
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L8").Value = 1420
$ws.Range("L9").Value = 123
$ws.Range("L10").Value = 141
$ws.Range("L11").Value = 353
$ws.Range("L13").Value = 31
$ws.Range("L15").Value = 187
$ws.Range("L20").Value = 541
$ws.Range("L22").Value = 72
$ws.Range("L27").Value = 186
$ws.Range("L29").Value = 1203
$ws.Range("L30").Value = 95
$ws.Range("L33").Value = 966
$ws.Range("L41").Value = 92
$ws.Range("L42").Value = 679
$ws.Range("L44").Value = 148
$ws.Range("L46").Value = 52
$ws.Range("L47").Value = 150
$ws.Range("L48").Value = 279
$ws.Range("L59").Value = 37
$ws.Range("L63").Value = 70
$ws.Range("L65").Value = 431
$ws.Range("L75").Value = 79
$ws.Range("L76").Value = 341
$ws.Range("L84").Value = 206
$ws.Range("L85").Value = 1072
$ws.Range("L89").Value = 289
$ws.Range("L90").Value = 233
$ws.Range("L95").Value = 300
$ws.Range("L97").Value = 172
$ws.Range("L99").Value = 373
$ws.Range("L100").Value = 42
$ws.Range("L101").Value = 21537
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 83
$ws.Range("L6").Value = 58
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 68
$ws.Range("L7").Value = 206
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 464
$ws.Range("L4").Value = 66
$ws.Range("L6").Value = 287
$ws.Range("L7").Value = 1203
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 113
$ws.Range("L7").Value = 279
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 58
$ws.Range("L7").Value = 148
$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 68
$ws.Range("L7").Value = 341
$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 92
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 235
$ws.Range("L6").Value = 193
$ws.Range("L7").Value = 679
$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("L5").Value = 16
$ws.Range("L6").Value = 31
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 141
$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L2").Value = 18
$ws.Range("L7").Value = 52
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 171
$ws.Range("L7").Value = 541
$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 42
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 150
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L3").Value = 60
$ws.Range("L7").Value = 187
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 99
$ws.Range("L7").Value = 353
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L3").Value = 46
$ws.Range("L7").Value = 123
$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 37
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L2").Value = 45
$ws.Range("L7").Value = 172
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 76
$ws.Range("L7").Value = 289
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L3").Value = 52
$ws.Range("L5").Value = 3
$ws.Range("L7").Value = 186
$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L4").Value = 9
$ws.Range("L7").Value = 79
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L3").Value = 68
$ws.Range("L7").Value = 233
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 319
$ws.Range("L7").Value = 1072
$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 72
